# Apply updated cryptocurrency price/volume data to Sheet1 (Aug 5 2024 refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "50.049.95"
$ws.Range("E2").Value = "  -17.70%  "

# Row 3
$ws.Range("D3").Value = "2.222.61"
$ws.Range("E3").Value = "  -23.53%  "

# Row 4
$ws.Range("E4").Value = "  +0.15%  "

# Row 5
$ws.Range("D5").Value = "'422.63"
$ws.Range("E5").Value = "  -19.61%  "

# Row 6
$ws.Range("D6").Value = "'112.74"
$ws.Range("E6").Value = "  -21.66%  "

# Row 7
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").Value = "'0.443"
$ws.Range("E8").Value = "  -18.65%  "

# Row 9
$ws.Range("D9").Value = "2.220.16"
$ws.Range("E9").Value = "  -23.80%  "

# Row 10
$ws.Range("D10").Value = "'5.02"
$ws.Range("E10").Value = "  -17.37%  "

# Row 11
$ws.Range("D11").Value = "'0.0833"
$ws.Range("E11").Value = "  -22.38%  "

# Row 12
$ws.Range("D12").Value = "'0.287"
$ws.Range("E12").Value = "  -19.84%  "

# Row 13
$ws.Range("D13").Value = "'0.118"
$ws.Range("E13").Value = "  -7.85%  "

# Row 14
$ws.Range("D14").Value = "2.614.10"
$ws.Range("E14").Value = "  -23.45%  "

# Row 15
$ws.Range("D15").Value = "50.351.47"
$ws.Range("E15").Value = "  -17.13%  "

# Row 16
$ws.Range("D16").Value = "'17.93"
$ws.Range("E16").Value = "  -20.41%  "

# Row 17
$ws.Range("D17").Value = "'0.0000111"
$ws.Range("E17").Value = "  -21.45%  "

# Row 18
$ws.Range("D18").Value = "2.240.29"
$ws.Range("E18").Value = "  -23.10%  "

# Row 19
$ws.Range("D19").Value = "'3.82"
$ws.Range("E19").Value = "  -21.79%  "

# Row 20
$ws.Range("D20").Value = "'286.06"
$ws.Range("E20").Value = "  -18.73%  "

# Row 21
$ws.Range("D21").Value = "'0.990"
$ws.Range("E21").Value = "  -0.95%  "

# Row 22
$ws.Range("D22").Value = "'5.62"
$ws.Range("E22").Value = "  -1.81%  "

# Row 23
$ws.Range("D23").Value = "'8.30"
$ws.Range("E23").Value = "  -28.24%  "

# Row 24
$ws.Range("D24").Value = "'4.89"
$ws.Range("E24").Value = "  -24.86%  "

# Row 25
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.09%  "

# Row 26
$ws.Range("D26").Value = "'51.89"
$ws.Range("E26").Value = "  -19.90%  "

# Row 27
$ws.Range("D27").Value = "'0.354"
$ws.Range("E27").Value = "  -21.35%  "

# Row 28
$ws.Range("D28").Value = "2.344.91"
$ws.Range("E28").Value = "  -22.67%  "

# Row 29
$ws.Range("D29").Value = "'0.131"
$ws.Range("E29").Value = "  -25.63%  "

# Row 30
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  -0.15%  "

# Row 31
$ws.Range("D31").Value = "'6.48"
$ws.Range("E31").Value = "  -17.21%  "

# Row 32
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").Value = "0.0₃0611"
$ws.Range("E32").Value = "  -29.08%  "

# Row 33
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "'139.17"
$ws.Range("E33").Value = "  -9.21%  "

# Row 34
$ws.Range("D34").Value = "'16.06"
$ws.Range("E34").Value = "  -18.11%  "

# Row 35
$ws.Range("D35").Value = "'1.27"
$ws.Range("E35").Value = "  -23.99%  "

# Row 36
$ws.Range("D36").Value = "'4.50"
$ws.Range("E36").Value = "  -18.86%  "

# Row 37
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  +0.22%  "

# Row 38
$ws.Range("D38").Value = "'3.21"
$ws.Range("E38").Value = "  -26.66%  "

# Row 39
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "'31.35"
$ws.Range("E39").Value = "  -16.53%  "

# Row 40
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "'0.935"
$ws.Range("E40").Value = "  -21.80%  "

# Row 41
$ws.Range("D41").Value = "'0.722"
$ws.Range("E41").Value = "  -27.25%  "

# Row 42
$ws.Range("D42").Value = "'10.12"
$ws.Range("E42").Value = "  -2.02%  "

# Row 43
$ws.Range("D43").Value = "'0.538"
$ws.Range("E43").Value = "  -17.53%  "

# Row 44
$ws.Range("D44").Value = "'2.99"
$ws.Range("E44").Value = "  -19.37%  "

# Row 45
$ws.Range("D45").Value = "'0.0477"
$ws.Range("E45").Value = "  -17.78%  "

# Row 46
$ws.Range("D46").Value = "1.821.13"
$ws.Range("E46").Value = "  -20.30%  "

# Row 47
$ws.Range("D47").Value = "'1.07"
$ws.Range("E47").Value = "  -27.10%  "

# Row 48
$ws.Range("D48").Value = "'0.0193"
$ws.Range("E48").Value = "  -18.23%  "

# Row 49
$ws.Range("D49").Value = "'0.0768"
$ws.Range("E49").Value = "  -16.22%  "

# Row 50
$ws.Range("B50").Value = "ZEEBU"
$ws.Range("C50").Value = "https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu"
$ws.Range("D50").Value = "'4.61"
$ws.Range("E50").Value = "  -5.67%  "

# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'15.08"
$ws.Range("E51").Value = "  -25.77%  "
